# Tasks_Tracking.xlsx update:
#  - Row 4 ("Backend Directories Creation( Packaging)") moves from
#    "in progress" (yellow) to "COMPLETED" (green), and its STATUS cell
#    (C4) is filled in with the text "COMPLETED".
#  - Row 5 ("SERVICE Layer Design") moves from "not started" (no fill)
#    to "in progress" (yellow).
#  - The active selection moves from E9 to F7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$green  = 5287936   # RGB(0, 176, 80)  -> matches existing fill "FF00B050"
$yellow = 65535      # RGB(255, 255, 0) -> matches existing fill "FFFFFF00"

# Row 4: mark as COMPLETED (green fill + status text)
$ws.Range("A4:C4").Interior.Color = $green
$ws.Range("C4").Value = "COMPLETED"

# Row 5: mark as in progress (yellow fill)
$ws.Range("A5:C5").Interior.Color = $yellow

# Update the saved selection to F7
$ws.Range("F7").Select()
